# edit.ps1 - apply updated cryptocurrency market data (prices / 1h volume %) to cryptos.xlsx
# Data pulled by the scheduled "Updated cryptos list" GitHub Actions job.
#
# Notes:
#  - Columns D (Price) and E (Volume(1h)) are stored as plain text in the sheet (values such as
#    "68.773.63" use dots as thousands separators and are not valid numbers, and the Volume(1h)
#    column keeps its literal padding spaces), so we must make sure Excel keeps every updated
#    value as text instead of silently re-interpreting it as a number.
#  - For the handful of new Price values that *do* look like plain numbers (e.g. "1.00", "576.86")
#    a leading single-quote is used - exactly like typing `'576.86` into a cell in Excel - which
#    forces a text entry while leaving the displayed/stored text identical to the source value.
#  - Rows 33 and 34 swapped rank (Mantle and NEARProtocol traded places), so columns B/C/D/E are
#    rewritten for both rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.647.27'
$ws.Range("E2").Value = '  -4.42%  '
$ws.Range("D3").Value = '3.504.52'
$ws.Range("E3").Value = '  -4.27%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").Value = '''576.86'
$ws.Range("E5").Value = '  -2.48%  '
$ws.Range("D6").Value = '''174.11'
$ws.Range("E6").Value = '  -3.90%  '
$ws.Range("E7").Value = '  -0.54%  '
$ws.Range("D8").Value = '3.495.50'
$ws.Range("E8").Value = '  -4.35%  '
$ws.Range("E9").Value = '  +0.13%  '
$ws.Range("E10").Value = '  -7.38%  '
$ws.Range("D11").Value = '''6.61'
$ws.Range("E11").Value = '  +7.38%  '
$ws.Range("D12").Value = '''0.599'
$ws.Range("E12").Value = '  -1.84%  '
$ws.Range("D13").Value = '''47.16'
$ws.Range("E13").Value = '  -5.64%  '
$ws.Range("E14").Value = '  -4.25%  '
$ws.Range("D15").Value = '''680.05'
$ws.Range("E15").Value = '  -0.37%  '
$ws.Range("D16").Value = '4.074.67'
$ws.Range("E16").Value = '  -3.89%  '
$ws.Range("D17").Value = '''8.85'
$ws.Range("E17").Value = '  -2.48%  '
$ws.Range("D18").Value = '68.818.47'
$ws.Range("E18").Value = '  -4.31%  '
$ws.Range("D19").Value = '3.503.90'
$ws.Range("E19").Value = '  -4.49%  '
$ws.Range("E20").Value = '  -1.47%  '
$ws.Range("E21").Value = '  -4.37%  '
$ws.Range("D22").Value = '''11.19'
$ws.Range("E22").Value = '  -4.18%  '
$ws.Range("E23").Value = '  -4.30%  '
$ws.Range("D24").Value = '''16.34'
$ws.Range("E24").Value = '  -8.65%  '
$ws.Range("D25").Value = '''97.45'
$ws.Range("E25").Value = '  -6.01%  '
$ws.Range("E26").Value = '  -5.24%  '
$ws.Range("D27").Value = '''1.00'
$ws.Range("E27").Value = '  -0.06%  '
$ws.Range("E28").Value = '  -6.64%  '
$ws.Range("D29").Value = '''9.38'
$ws.Range("E29").Value = '  -8.17%  '
$ws.Range("D30").Value = '''32.96'
$ws.Range("E30").Value = '  -6.90%  '
$ws.Range("E31").Value = '  -5.08%  '
$ws.Range("E32").Value = '  -9.15%  '
$ws.Range("B33").Value = 'NEARProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D33").Value = '''7.27'
$ws.Range("E33").Value = '  -1.22%  '
$ws.Range("B34").Value = 'Mantle'
$ws.Range("C34").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D34").Value = '''1.35'
$ws.Range("E34").Value = '  -6.30%  '
$ws.Range("D35").Value = '''565.71'
$ws.Range("E35").Value = '  -2.30%  '
$ws.Range("D36").Value = '''3.62'
$ws.Range("E36").Value = '  -13.91%  '
$ws.Range("D37").Value = '''10.85'
$ws.Range("E37").Value = '  -4.41%  '
$ws.Range("E38").Value = '  -3.59%  '
$ws.Range("D39").Value = '''57.08'
$ws.Range("E39").Value = '  -4.17%  '
$ws.Range("D40").Value = '''0.998'
$ws.Range("E40").Value = '  -0.15%  '
$ws.Range("D41").Value = '''0.137'
$ws.Range("E41").Value = '  -4.54%  '
$ws.Range("E42").Value = '  -6.00%  '
$ws.Range("D43").Value = '3.449.61'
$ws.Range("E43").Value = '  -7.63%  '
$ws.Range("D44").Value = '''0.335'
$ws.Range("E44").Value = '  -3.50%  '
$ws.Range("D45").Value = '''33.42'
$ws.Range("E45").Value = '  -6.47%  '
$ws.Range("D46").Value = '0.0₃0699'
$ws.Range("E46").Value = '  -8.43%  '
$ws.Range("D47").Value = '''2.90'
$ws.Range("E47").Value = '  +2.80%  '
$ws.Range("D48").Value = '''2.59'
$ws.Range("E48").Value = '  -7.71%  '
$ws.Range("E49").Value = '  -0.67%  '
$ws.Range("D50").Value = '''134.11'
$ws.Range("E50").Value = '  +0.98%  '
$ws.Range("E51").Value = '  -0.97%  '
